$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.650.27"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.602.00"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.77"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.51%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "1.831.24"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "1.603.90"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("D15").Value = "29.643.66"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.94"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.96"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("E19").Value = "  +3.85%  "
$ws.Range("D20").Value = "0.0₃0697"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.40"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.20"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.43"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  +2.64%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("D34").Value = "1.424.97"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +3.41%  "
$ws.Range("E36").Value = "  +4.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "58.55"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.545"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.49%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0499"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.85%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.97"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.816"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.32"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.975"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +16.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").Value = "1.741.70"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.74"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("D51").Value = "0.0₆0103"
$ws.Range("E51").Value = "  +2.96%  "
